$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45189 = 2023-09-20)
# that was bumped by one day (45190 = 2023-09-21) for every data row.
$ws.Range("C2:C246").Value = 45190
